$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# All_SANs: append 3 new rows (147-149) for new "Laptop 840 G10" SAN entries
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All_SANs")

$wsAll.Cells.Item(147, 1).Value = "SAN457865"
$wsAll.Cells.Item(147, 2).Value = "Laptop 840 G10"
$wsAll.Cells.Item(147, 3).Value = "2024-11-17 20:05:59"
$wsAll.Cells.Item(147, 4).Value = "'4.2"

$wsAll.Cells.Item(148, 1).Value = "SAN125335"
$wsAll.Cells.Item(148, 2).Value = "Laptop 840 G10"
$wsAll.Cells.Item(148, 3).Value = "2024-11-17 23:27:40"
$wsAll.Cells.Item(148, 4).Value = "'4.2"

$wsAll.Cells.Item(149, 1).Value = "SAN125448"
$wsAll.Cells.Item(149, 2).Value = "Laptop 840 G10"
$wsAll.Cells.Item(149, 3).Value = "2024-11-17 23:27:42"
$wsAll.Cells.Item(149, 4).Value = "'4.2"

# The source file has no explicit cell-level styles on these rows, but the
# sheet's <cols> definition applies a default style to columns A-C; resetting
# to "Normal" here matches the unstyled cells produced upstream.
$wsAll.Range("A147:D149").Style = "Normal"

# ---------------------------------------------------------------------------
# 4.2_Items: bump "Laptop 840 G10" LastCount / NewCount totals (row 6)
# ---------------------------------------------------------------------------
$ws42Items = $wb.Worksheets.Item("4.2_Items")
$ws42Items.Cells.Item(6, 2).Value = 39
$ws42Items.Cells.Item(6, 3).Value = 41

# ---------------------------------------------------------------------------
# 4.2_Timestamps: fill in the previously-blank placeholder rows 44-48
# ---------------------------------------------------------------------------
$ws42Ts = $wb.Worksheets.Item("4.2_Timestamps")

$ws42Ts.Cells.Item(44, 1).Value = "2024-11-17 20:05:59"
$ws42Ts.Cells.Item(44, 2).Value = "Laptop 840 G10"
$ws42Ts.Cells.Item(44, 3).Value = "add"
$ws42Ts.Cells.Item(44, 4).Value = "SAN457865"

$ws42Ts.Cells.Item(45, 1).Value = "2024-11-17 20:05:59"
$ws42Ts.Cells.Item(45, 2).Value = "Laptop 840 G10"
$ws42Ts.Cells.Item(45, 3).Value = "add 1"

$ws42Ts.Cells.Item(46, 1).Value = "2024-11-17 23:27:40"
$ws42Ts.Cells.Item(46, 2).Value = "Laptop 840 G10"
$ws42Ts.Cells.Item(46, 3).Value = "add"
$ws42Ts.Cells.Item(46, 4).Value = "SAN125335"

$ws42Ts.Cells.Item(47, 1).Value = "2024-11-17 23:27:42"
$ws42Ts.Cells.Item(47, 2).Value = "Laptop 840 G10"
$ws42Ts.Cells.Item(47, 3).Value = "add"
$ws42Ts.Cells.Item(47, 4).Value = "SAN125448"

$ws42Ts.Cells.Item(48, 1).Value = "2024-11-17 23:27:43"
$ws42Ts.Cells.Item(48, 2).Value = "Laptop 840 G10"
$ws42Ts.Cells.Item(48, 3).Value = "add 2"

# ---------------------------------------------------------------------------
# B4.3_Items: bump "Laptop 840 G6" LastCount / NewCount totals (row 2)
# ---------------------------------------------------------------------------
$wsB43Items = $wb.Worksheets.Item("B4.3_Items")
$wsB43Items.Cells.Item(2, 2).Value = 2
$wsB43Items.Cells.Item(2, 3).Value = 4

# ---------------------------------------------------------------------------
# B4.3_Timestamps: append new log row 4
# ---------------------------------------------------------------------------
$wsB43Ts = $wb.Worksheets.Item("B4.3_Timestamps")
$wsB43Ts.Cells.Item(4, 1).Value = "2024-11-17 23:48:42"
$wsB43Ts.Cells.Item(4, 2).Value = "Laptop 840 G6"
$wsB43Ts.Cells.Item(4, 3).Value = "add 2"

# ---------------------------------------------------------------------------
# Darwin_Timestamps: drop the stray empty D34 cell
# ---------------------------------------------------------------------------
$wsDarwinTs = $wb.Worksheets.Item("Darwin_Timestamps")
$wsDarwinTs.Cells.Item(34, 4).ClearContents()
